# "taking costs out of infrastructure database"
#
# The RESOURCES sheet had a "costs_kWh" column (column E, between CO2 and
# reference) holding operating-cost figures (with a cell comment explaining
# the units). This column is removed entirely: the reference column (F)
# shifts left into its place, the now-unused "costs_kWh" shared string is
# dropped, and the comment that was anchored on E1 goes away with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RESOURCES")

# Remove the cell comment living on the costs column header before the
# column shifts, so it doesn't survive anchored to the wrong cell.
$costsComment = $ws.Range("E1").Comment
if ($costsComment -ne $null) {
    $costsComment.Delete()
}

# Delete the whole costs_kWh column; everything to the right (the
# "reference" column) shifts left to take its place.
$ws.Columns("E:E").Delete()

# Match the author's final selection/cursor position on the sheet.
$ws.Range("F9").Select()
